$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 4
$ws.Range("B4").Value = "Jerrygood"
$ws.Range("C4").Value = "deepakverma.knp2019@gmail.com"
$ws.Range("D4").Value = "xyz"
$ws.Range("E4").Value = "xyz"
